$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update existing rows 20-31 with the cells that changed (data shifted down by one,
#     row 20 becomes a brand-new record, row 21 ["Segunda"] stays put) ---
# Row 20
$ws.Range("D20").Value = 44818
$ws.Range("J20").Value = 60
$ws.Range("K20").Value = 15000
$ws.Range("L20").Value = 15000
$ws.Range("M20").Value = 15000
$ws.Range("P20").Value = 833

# Row 22
$ws.Range("D22").Value = 44384
$ws.Range("I22").Value = "Segunda"
$ws.Range("J22").Value = 60
$ws.Range("K22").Value = 15000
$ws.Range("M22").Value = 15000
$ws.Range("P22").Value = 833

# Row 23
$ws.Range("D23").Value = 44813
$ws.Range("J23").Value = 100
$ws.Range("K23").Value = 14000
$ws.Range("L23").Value = 15000
$ws.Range("M23").Value = 14500
$ws.Range("P23").Value = 806

# Row 24
$ws.Range("D24").Value = 44664
$ws.Range("J24").Value = 160

# Row 25
$ws.Range("D25").Value = 44651
$ws.Range("J25").Value = 60

# Row 26
$ws.Range("D26").Value = 44656
$ws.Range("K26").Value = 15000
$ws.Range("L26").Value = 16000
$ws.Range("M26").Value = 15500
$ws.Range("P26").Value = 861

# Row 27
$ws.Range("D27").Value = 44775
$ws.Range("J27").Value = 100
$ws.Range("K27").Value = 17000
$ws.Range("L27").Value = 18000
$ws.Range("M27").Value = 17500
$ws.Range("P27").Value = 972

# Row 28
$ws.Range("D28").Value = 44809
$ws.Range("J28").Value = 60
$ws.Range("K28").Value = 14000
$ws.Range("L28").Value = 15000
$ws.Range("M28").Value = 14500
$ws.Range("P28").Value = 806

# Row 29
$ws.Range("D29").Value = 44761
$ws.Range("J29").Value = 100

# Row 30
$ws.Range("D30").Value = 44763
$ws.Range("J30").Value = 80

# Row 31
$ws.Range("D31").Value = 44771
$ws.Range("J31").Value = 60

# --- Append new row 32 (duplicate of what used to be row 31) ---
$ws.Range("D32").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("A32").Value = 7
$ws.Range("B32").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C32").Value = "Ñuble"
$ws.Range("D32").Value = 44782
$ws.Range("E32").Value = 16
$ws.Range("F32").Value = 100112043
$ws.Range("G32").Value = "Pepino dulce"
$ws.Range("H32").Value = "Cultivar IV Región"
$ws.Range("I32").Value = "Primera"
$ws.Range("J32").Value = 120
$ws.Range("K32").Value = 17000
$ws.Range("L32").Value = 18000
$ws.Range("M32").Value = 17500
$ws.Range("N32").Value = "$/bandeja 18 kilos"
$ws.Range("O32").Value = "Provincia de Limarí"
$ws.Range("P32").Value = 972
$ws.Range("Q32").Value = 18
$ws.Range("R32").Value = "Hortaliza"
